$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width change (39.5546875 -> 58) ---
# Excel's ColumnWidth property is expressed in "characters" and gets
# padded/rounded to pixels internally before being stored back as the
# raw OOXML column width units, so the COM value that round-trips to a
# stored width of exactly 58 is 400/7 (~57.142857) characters.
$ws.Columns.Item(2).ColumnWidth = 400/7

# --- Row 19 updates: D19 date, F19 status ---
$ws.Range("D19").Value = 44897
$ws.Range("F19").Value = "Done"

# --- Row 20 updates: D20 date, F20 status ---
$ws.Range("D20").Value = 44897
$ws.Range("F20").Value = "Done"

# --- Row 21 update: D21 date ---
$ws.Range("D21").Value = 44898

# --- New underlined-font formatting marks on F23 and E25 ---
$ws.Range("F23").Font.Underline = $true
$ws.Range("E25").Font.Underline = $true

# --- Selection moved from A22 to E25 ---
[void]$ws.Range("E25").Select()
